# Update NATMI LR-pair (Gpi1-Amfr) expression/specificity values with new TPM-based results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 50.000707
$ws.Range("H2").Value = 150.002121
$ws.Range("I2").Value = 0.3101889378151254
$ws.Range("J2").Value = 0.3101889378151254
$ws.Range("M2").Value = 17.051494
$ws.Range("N2").Value = 51.154482
$ws.Range("O2").Value = 0.2221091776636934
$ws.Range("P2").Value = 0.2221091776636934
$ws.Range("Q2").Value = 852.586755406258
$ws.Range("R2").Value = 7673.280798656322
$ws.Range("S2").Value = 0.06889580989849203
$ws.Range("T2").Value = 0.06889580989849201

# Row 3
$ws.Range("G3").Value = 50.000707
$ws.Range("H3").Value = 150.002121
$ws.Range("I3").Value = 0.3101889378151254
$ws.Range("J3").Value = 0.3101889378151254
$ws.Range("N3").Value = 64.372913
$ws.Range("O3").Value = 0.2795026791639973
$ws.Range("P3").Value = 0.2795026791639973
$ws.Range("Q3").Value = 1072.897053883164
$ws.Range("R3").Value = 9656.073484948472
$ws.Range("S3").Value = 0.08669863916636213
$ws.Range("T3").Value = 0.08669863916636208

# Row 4
$ws.Range("G4").Value = 50.000707
$ws.Range("H4").Value = 150.002121
$ws.Range("I4").Value = 0.3101889378151254
$ws.Range("J4").Value = 0.3101889378151254
$ws.Range("M4").Value = 21.25262066666667
$ws.Range("N4").Value = 63.757862
$ws.Range("O4").Value = 0.2768321708040216
$ws.Range("P4").Value = 0.2768321708040215
$ws.Range("Q4").Value = 1062.646058936145
$ws.Range("R4").Value = 9563.814530425301
$ws.Range("S4").Value = 0.08587027701475482
$ws.Range("T4").Value = 0.08587027701475479

# Row 5
$ws.Range("G5").Value = 50.000707
$ws.Range("H5").Value = 150.002121
$ws.Range("I5").Value = 0.3101889378151254
$ws.Range("J5").Value = 0.3101889378151254
$ws.Range("M5").Value = 17.009024
$ws.Range("N5").Value = 51.027072
$ws.Range("O5").Value = 0.2215559723682878
$ws.Range("P5").Value = 0.2215559723682878
$ws.Range("Q5").Value = 850.463225379968
$ws.Range("R5").Value = 7654.169028419712
$ws.Range("S5").Value = 0.06872421173551646
$ws.Range("T5").Value = 0.06872421173551645

# Row 6
$ws.Range("I6").Value = 0.2298958220161207
$ws.Range("J6").Value = 0.2298958220161207
$ws.Range("M6").Value = 17.051494
$ws.Range("N6").Value = 51.154482
$ws.Range("O6").Value = 0.2221091776636934
$ws.Range("P6").Value = 0.2221091776636934
$ws.Range("Q6").Value = 631.8927243337088
$ws.Range("R6").Value = 5687.034519003378
$ws.Range("S6").Value = 0.0510619719763194
$ws.Range("T6").Value = 0.05106197197631938

# Row 7
$ws.Range("I7").Value = 0.2298958220161207
$ws.Range("J7").Value = 0.2298958220161207
$ws.Range("N7").Value = 64.372913
$ws.Range("O7").Value = 0.2795026791639973
$ws.Range("P7").Value = 0.2795026791639973
$ws.Range("Q7").Value = 795.1751983113975
$ws.Range("R7").Value = 7156.576784802576
$ws.Range("S7").Value = 0.06425649818211523
$ws.Range("T7").Value = 0.06425649818211521

# Row 8
$ws.Range("I8").Value = 0.2298958220161207
$ws.Range("J8").Value = 0.2298958220161207
$ws.Range("M8").Value = 21.25262066666667
$ws.Range("N8").Value = 63.757862
$ws.Range("O8").Value = 0.2768321708040216
$ws.Range("P8").Value = 0.2768321708040215
$ws.Range("Q8").Value = 787.5776968452666
$ws.Range("R8").Value = 7088.199271607398
$ws.Range("S8").Value = 0.06364255946749768
$ws.Range("T8").Value = 0.06364255946749765

# Row 9
$ws.Range("I9").Value = 0.2298958220161207
$ws.Range("J9").Value = 0.2298958220161207
$ws.Range("M9").Value = 17.009024
$ws.Range("N9").Value = 51.027072
$ws.Range("O9").Value = 0.2215559723682878
$ws.Range("P9").Value = 0.2215559723682878
$ws.Range("Q9").Value = 630.3188749101654
$ws.Range("R9").Value = 5672.869874191488
$ws.Range("S9").Value = 0.05093479239018845
$ws.Range("T9").Value = 0.05093479239018844

# Row 10
$ws.Range("G10").Value = 46.58340866666666
$ws.Range("H10").Value = 139.750226
$ws.Range("I10").Value = 0.2889890747769074
$ws.Range("J10").Value = 0.2889890747769074
$ws.Range("M10").Value = 17.051494
$ws.Range("N10").Value = 51.154482
$ws.Range("O10").Value = 0.2221091776636934
$ws.Range("P10").Value = 0.2221091776636934
$ws.Range("Q10").Value = 794.3167133792147
$ws.Range("R10").Value = 7148.850420412932
$ws.Range("S10").Value = 0.06418712575249051
$ws.Range("T10").Value = 0.0641871257524905

# Row 11
$ws.Range("G11").Value = 46.58340866666666
$ws.Range("H11").Value = 139.750226
$ws.Range("I11").Value = 0.2889890747769074
$ws.Range("J11").Value = 0.2889890747769074
$ws.Range("N11").Value = 64.372913
$ws.Range("O11").Value = 0.2795026791639973
$ws.Range("P11").Value = 0.2795026791639973
$ws.Range("Q11").Value = 999.5699044475931
$ws.Range("R11").Value = 8996.129140028337
$ws.Range("S11").Value = 0.08077322064927039
$ws.Range("T11").Value = 0.08077322064927038

# Row 12
$ws.Range("G12").Value = 46.58340866666666
$ws.Range("H12").Value = 139.750226
$ws.Range("I12").Value = 0.2889890747769074
$ws.Range("J12").Value = 0.2889890747769074
$ws.Range("M12").Value = 21.25262066666667
$ws.Range("N12").Value = 63.757862
$ws.Range("O12").Value = 0.2768321708040216
$ws.Range("P12").Value = 0.2768321708040215
$ws.Range("Q12").Value = 990.0195137529792
$ws.Range("R12").Value = 8910.175623776813
$ws.Range("S12").Value = 0.08000147290913701
$ws.Range("T12").Value = 0.08000147290913698

# Row 13
$ws.Range("G13").Value = 46.58340866666666
$ws.Range("H13").Value = 139.750226
$ws.Range("I13").Value = 0.2889890747769074
$ws.Range("J13").Value = 0.2889890747769074
$ws.Range("M13").Value = 17.009024
$ws.Range("N13").Value = 51.027072
$ws.Range("O13").Value = 0.2215559723682878
$ws.Range("P13").Value = 0.2215559723682878
$ws.Range("Q13").Value = 792.3383160131413
$ws.Range("R13").Value = 7131.044844118273
$ws.Range("S13").Value = 0.06402725546600956
$ws.Range("T13").Value = 0.06402725546600956

# Row 14
$ws.Range("G14").Value = 27.55233366666667
$ws.Range("H14").Value = 82.65700100000001
$ws.Range("I14").Value = 0.1709261653918464
$ws.Range("J14").Value = 0.1709261653918464
$ws.Range("M14").Value = 17.051494
$ws.Range("N14").Value = 51.154482
$ws.Range("O14").Value = 0.2221091776636934
$ws.Range("P14").Value = 0.2221091776636934
$ws.Range("Q14").Value = 469.8084522031647
$ws.Range("R14").Value = 4228.276069828483
$ws.Range("S14").Value = 0.03796427003639146
$ws.Range("T14").Value = 0.03796427003639145

# Row 15
$ws.Range("G15").Value = 27.55233366666667
$ws.Range("H15").Value = 82.65700100000001
$ws.Range("I15").Value = 0.1709261653918464
$ws.Range("J15").Value = 0.1709261653918464
$ws.Range("N15").Value = 64.372913
$ws.Range("O15").Value = 0.2795026791639973
$ws.Range("P15").Value = 0.2795026791639973
$ws.Range("Q15").Value = 591.2079926904348
$ws.Range("R15").Value = 5320.871934213913
$ws.Range("S15").Value = 0.0477743211662496
$ws.Range("T15").Value = 0.04777432116624959

# Row 16
$ws.Range("G16").Value = 27.55233366666667
$ws.Range("H16").Value = 82.65700100000001
$ws.Range("I16").Value = 0.1709261653918464
$ws.Range("J16").Value = 0.1709261653918464
$ws.Range("M16").Value = 21.25262066666667
$ws.Range("N16").Value = 63.757862
$ws.Range("O16").Value = 0.2768321708040216
$ws.Range("P16").Value = 0.2768321708040215
$ws.Range("Q16").Value = 585.5592958990959
$ws.Range("R16").Value = 5270.033663091863
$ws.Range("S16").Value = 0.04731786141263207
$ws.Range("T16").Value = 0.04731786141263206

# Row 17
$ws.Range("G17").Value = 27.55233366666667
$ws.Range("H17").Value = 82.65700100000001
$ws.Range("I17").Value = 0.1709261653918464
$ws.Range("J17").Value = 0.1709261653918464
$ws.Range("M17").Value = 17.009024
$ws.Range("N17").Value = 51.027072
$ws.Range("O17").Value = 0.2215559723682878
$ws.Range("P17").Value = 0.2215559723682878
$ws.Range("Q17").Value = 468.6383045923414
$ws.Range("R17").Value = 4217.744741331073
$ws.Range("S17").Value = 0.03786971277657331
$ws.Range("T17").Value = 0.03786971277657331
